$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.387.35'
$ws.Range("E2").Value = '  -1.17%  '

$ws.Range("D3").Value = '3.071.60'
$ws.Range("E3").Value = '  -2.44%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.70%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.547'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.10%  '

$ws.Range("D9").Value = '3.075.90'
$ws.Range("E9").Value = '  -2.05%  '

$ws.Range("E10").Value = '  -3.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.87'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.42%  '

$ws.Range("E12").Value = '  +0.71%  '

$ws.Range("E13").Value = '  -1.94%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.31'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.26%  '

$ws.Range("E15").Value = '  -1.66%  '

$ws.Range("D16").Value = '3.580.39'
$ws.Range("E16").Value = '  -2.41%  '

$ws.Range("E17").Value = '  -0.37%  '

$ws.Range("D18").Value = '63.409.72'
$ws.Range("E18").Value = '  -0.84%  '

$ws.Range("D19").Value = '3.078.11'
$ws.Range("E19").Value = '  -2.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '475.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.719'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.64%  '

$ws.Range("E23").Value = '  +0.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.59%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.10'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.57%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.99'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.50%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.40%  '

$ws.Range("E30").Value = '  -1.04%  '

$ws.Range("B31").Value = 'FirstDigitalUSD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.24%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.37%  '

$ws.Range("E33").Value = '  +4.24%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.74%  '

$ws.Range("E35").Value = '  +2.82%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.34%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.59%  '

$ws.Range("E39").Value = '  -3.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.33'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.47'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.66%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '447.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.59%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.286'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.46%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.63'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.78%  '

$ws.Range("E45").Value = '  -2.26%  '

$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.111'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.38%  '

$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.816.61'
$ws.Range("E47").Value = '  -3.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.15'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.30'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.69%  '

$ws.Range("E51").Value = '  +0.45%  '

